$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  ,@('2025-11-29T23:42:12.152Z', 'HTTP', 'GET', '/admin/data', 200, 'OK', '', '', '', '', '::1')
  ,@('2025-11-30T02:11:31.838Z', 'HTTP', 'GET', '/user/6BF02F00?roomID=106', 404, 'ERROR', '', '', '106', 'UID 6BF02F00 not found in /user', '::ffff:172.28.219.204')
  ,@('2025-11-30T02:11:32.004Z', 'HTTP', 'GET', '/uid-name/6BF02F00', 404, 'ERROR', '', '', '', 'UID 6BF02F00 not found in /uid-name', '::ffff:172.28.219.204')
  ,@('2025-11-30T02:11:42.002Z', 'HTTP', 'GET', '/favicon.ico', 404, 'ERROR', '', '', '', '', '::1')
  ,@('2025-11-30T02:11:45.608Z', 'HTTP', 'GET', '/admin/data', 200, 'OK', '', '', '', '', '::1')
  ,@('2025-11-30T02:11:47.226Z', 'HTTP', 'POST', '/register', 200, 'OK', 'YAHYA', '6BF02F00', '106', 'Registration successful for new UID', '::ffff:172.28.219.204')
  ,@('2025-11-30T02:11:50.060Z', 'HTTP', 'GET', '/user/6BF02F00?roomID=106', 200, 'OK', 'YAHYA', '6BF02F00', '106', 'Access granted and counter incremented', '::ffff:172.28.219.204')
  ,@('2025-11-30T02:11:50.151Z', 'HTTP', 'GET', '/uid-name/6BF02F00', 200, 'OK', 'YAHYA', '6BF02F00', '', 'UID to username lookup success', '::ffff:172.28.219.204')
  ,@('2025-11-30T02:11:50.557Z', 'HTTP', 'GET', '/admin/data', 200, 'OK', '', '', '', '', '::1')
  ,@('2025-11-30T02:11:52.604Z', 'HTTP', 'GET', '/admin/data', 200, 'OK', '', '', '', '', '::1')
  ,@('2025-11-30T02:11:55.935Z', 'HTTP', 'GET', '/user/6BF02F00?roomID=106', 200, 'OK', 'YAHYA', '6BF02F00', '106', 'Access granted and counter incremented', '::ffff:172.28.219.204')
  ,@('2025-11-30T02:11:56.012Z', 'HTTP', 'GET', '/uid-name/6BF02F00', 200, 'OK', 'YAHYA', '6BF02F00', '', 'UID to username lookup success', '::ffff:172.28.219.204')
  ,@('2025-11-30T02:11:56.750Z', 'HTTP', 'GET', '/admin/data', 200, 'OK', '', '', '', '', '::1')
  ,@('2025-11-30T02:11:59.482Z', 'HTTP', 'GET', '/user/6BF02F00?roomID=106', 200, 'OK', 'YAHYA', '6BF02F00', '106', 'Access granted and counter incremented', '::ffff:172.28.219.204')
  ,@('2025-11-30T02:11:59.589Z', 'HTTP', 'GET', '/uid-name/6BF02F00', 200, 'OK', 'YAHYA', '6BF02F00', '', 'UID to username lookup success', '::ffff:172.28.219.204')
  ,@('2025-11-30T02:12:00.828Z', 'HTTP', 'GET', '/user/6BF02F00?roomID=106', 200, 'OK', 'YAHYA', '6BF02F00', '106', 'Access granted and counter incremented', '::ffff:172.28.219.204')
  ,@('2025-11-30T02:12:00.897Z', 'HTTP', 'GET', '/uid-name/6BF02F00', 200, 'OK', 'YAHYA', '6BF02F00', '', 'UID to username lookup success', '::ffff:172.28.219.204')
  ,@('2025-11-30T02:12:01.743Z', 'HTTP', 'GET', '/admin/data', 200, 'OK', '', '', '', '', '::1')
  ,@('2025-11-30T02:12:01.961Z', 'HTTP', 'GET', '/user/6BF02F00?roomID=106', 200, 'OK', 'YAHYA', '6BF02F00', '106', 'Access granted and counter incremented', '::ffff:172.28.219.204')
  ,@('2025-11-30T02:12:02.064Z', 'HTTP', 'GET', '/uid-name/6BF02F00', 200, 'OK', 'YAHYA', '6BF02F00', '', 'UID to username lookup success', '::ffff:172.28.219.204')
  ,@('2025-11-30T02:12:06.750Z', 'HTTP', 'GET', '/admin/data', 200, 'OK', '', '', '', '', '::1')
  ,@('2025-11-30T02:12:07.739Z', 'HTTP', 'DELETE', '/room/6BF02F00/106', 200, 'OK', 'YAHYA', '6BF02F00', '106', 'User deleted because no rooms remain', '::ffff:172.28.219.204')
  ,@('2025-11-30T02:12:11.741Z', 'HTTP', 'GET', '/admin/data', 200, 'OK', '', '', '', '', '::1')
  ,@('2025-11-30T02:12:16.764Z', 'HTTP', 'GET', '/admin/data', 200, 'OK', '', '', '', '', '::1')
  ,@('2025-11-30T02:12:21.753Z', 'HTTP', 'GET', '/admin/data', 200, 'OK', '', '', '', '', '::1')
  ,@('2025-11-30T02:12:26.751Z', 'HTTP', 'GET', '/admin/data', 200, 'OK', '', '', '', '', '::1')
)

$startRow = 1386
for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $startRow + $i
  $rowData = $rows[$i]
  for ($c = 1; $c -le $rowData.Count; $c++) {
    $val = $rowData[$c - 1]
    $cell = $ws.Cells.Item($r, $c)
    if ($val -is [int]) {
      $cell.Value = $val
    } elseif ($val -eq "") {
      $cell.Value = "'"
    } elseif ($val -match '^[0-9]+$') {
      $cell.NumberFormat = "@"
      $cell.Value = $val
    } else {
      $cell.Value = $val
    }
  }
  $rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, $rowData.Count))
  $rowRange.ClearFormats()
}
